$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 23036508
$ws.Range("I11").Value = 23036508
$ws.Range("K11").Value = 23036508
$ws.Range("M11").Value = -23036368
$ws.Range("H15").Value = 86686.336
$ws.Range("I15").Value = 86686.336
$ws.Range("K15").Value = 260059.008
$ws.Range("M15").Value = -259890.008
$ws.Range("H17").Value = 2820.0645
$ws.Range("J17").Value = 2820.0645
$ws.Range("L17").Value = 8460.193499999999
$ws.Range("N17").Value = -8796.193499999999
$ws.Range("H28").Value = 794412.5600000001
$ws.Range("I28").Value = 1234810
$ws.Range("J28").Value = 1697.2
$ws.Range("K28").Value = 1234810
$ws.Range("L28").Value = 1697.2
$ws.Range("M28").Value = -1234325
$ws.Range("N28").Value = -2667.2
$ws.Range("H107").Value = 926468
$ws.Range("I107").Value = 1111611
$ws.Range("J107").Value = 753
$ws.Range("K107").Value = 1111611
$ws.Range("L107").Value = 753
$ws.Range("M107").Value = -1109691
$ws.Range("N107").Value = -4593
$ws.Range("H111").Value = 1793.6666
$ws.Range("I111").Value = 1410
$ws.Range("J111").Value = 2177.3333
$ws.Range("K111").Value = 4230
$ws.Range("L111").Value = 6531.999899999999
$ws.Range("M111").Value = -1163
$ws.Range("N111").Value = -12665.9999
$ws.Range("H112").Value = 11859011
$ws.Range("J112").Value = 12398011
$ws.Range("L112").Value = 37194033
$ws.Range("N112").Value = -37196249
$ws.Range("H115").Value = 2898.6365
$ws.Range("I115").Value = 735.625
$ws.Range("J115").Value = 8666.666999999999
$ws.Range("K115").Value = 2206.875
$ws.Range("L115").Value = 26000.001
$ws.Range("M115").Value = -639.875
$ws.Range("N115").Value = -29134.001
$ws.Range("H132").Value = 54691.43
$ws.Range("I132").Value = 58139.723
$ws.Range("K132").Value = 174419.169
$ws.Range("M132").Value = -171889.169
$ws.Range("H139").Value = 41563.332
$ws.Range("J139").Value = 41563.332
$ws.Range("L139").Value = 41563.332
$ws.Range("N139").Value = -51843.332
$ws.Range("H140").Value = 63743.332
$ws.Range("J140").Value = 63743.332
$ws.Range("L140").Value = 63743.332
$ws.Range("N140").Value = -74103.33199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 24240.418
$ws.Range("I2").Value = 33201.453
$ws.Range("J2").Value = 1091.0834
$ws.Range("K2").Value = 33201.453
$ws.Range("L2").Value = 1091.0834
$ws.Range("M2").Value = -33088.453
$ws.Range("N2").Value = -1317.0834
$ws.Range("H4").Value = 1192.7778
$ws.Range("I4").Value = 1192.7778
$ws.Range("K4").Value = 1192.7778
$ws.Range("M4").Value = -1076.7778
$ws.Range("H45").Value = 1716.6666
$ws.Range("I45").Value = 1600
$ws.Range("J45").Value = 1950
$ws.Range("K45").Value = 1600
$ws.Range("L45").Value = 1950
$ws.Range("M45").Value = -1223
$ws.Range("N45").Value = -2704
$ws.Range("H63").Value = 6606.4287
$ws.Range("I63").Value = 6457.5
$ws.Range("J63").Value = 7500
$ws.Range("K63").Value = 6457.5
$ws.Range("L63").Value = 7500
$ws.Range("M63").Value = -5771.5
$ws.Range("N63").Value = -8872
$ws.Range("H66").Value = 6606.4287
$ws.Range("I66").Value = 6457.5
$ws.Range("J66").Value = 7500
$ws.Range("K66").Value = 32287.5
$ws.Range("L66").Value = 37500
$ws.Range("M66").Value = -28855.5
$ws.Range("N66").Value = -44364
$ws.Range("H116").Value = 24240.418
$ws.Range("I116").Value = 33201.453
$ws.Range("J116").Value = 1091.0834
$ws.Range("K116").Value = 33201.453
$ws.Range("L116").Value = 1091.0834
$ws.Range("M116").Value = -30907.453
$ws.Range("N116").Value = -5679.0834
$ws.Range("H122").Value = 2778.027
$ws.Range("I122").Value = 1959.6897
$ws.Range("J122").Value = 5744.5
$ws.Range("K122").Value = 5879.0691
$ws.Range("L122").Value = 17233.5
$ws.Range("M122").Value = -3429.0691
$ws.Range("N122").Value = -22133.5
$ws.Range("H132").Value = 2010.0546
$ws.Range("I132").Value = 1734.8
$ws.Range("J132").Value = 3248.7
$ws.Range("K132").Value = 5204.4
$ws.Range("L132").Value = 9746.099999999999
$ws.Range("M132").Value = -2674.4
$ws.Range("N132").Value = -14806.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 24240.418
$ws.Range("I3").Value = 33201.453
$ws.Range("J3").Value = 1091.0834
$ws.Range("K3").Value = 33201.453
$ws.Range("L3").Value = 1091.0834
$ws.Range("M3").Value = -33087.453
$ws.Range("N3").Value = -1319.0834
$ws.Range("H37").Value = 16750
$ws.Range("I37").Value = 5666.6665
$ws.Range("K37").Value = 5666.6665
$ws.Range("M37").Value = -5529.6665
$ws.Range("H64").Value = 1937.5
$ws.Range("I64").Value = 2975
$ws.Range("J64").Value = 900
$ws.Range("K64").Value = 2975
$ws.Range("L64").Value = 900
$ws.Range("M64").Value = -2750
$ws.Range("N64").Value = -1350
$ws.Range("H67").Value = 1937.5
$ws.Range("I67").Value = 2975
$ws.Range("J67").Value = 900
$ws.Range("K67").Value = 2975
$ws.Range("L67").Value = 900
$ws.Range("M67").Value = -2195
$ws.Range("N67").Value = -2460
$ws.Range("H80").Value = 669.4815
$ws.Range("I80").Value = 257
$ws.Range("J80").Value = 787.3333
$ws.Range("K80").Value = 257
$ws.Range("L80").Value = 787.3333
$ws.Range("M80").Value = 741
$ws.Range("N80").Value = -2783.3333
$ws.Range("H83").Value = 669.4815
$ws.Range("I83").Value = 257
$ws.Range("J83").Value = 787.3333
$ws.Range("K83").Value = 1285
$ws.Range("L83").Value = 3936.6665
$ws.Range("M83").Value = 3707
$ws.Range("N83").Value = -13920.6665
$ws.Range("H86").Value = 11620.7
$ws.Range("I86").Value = 1966.6666
$ws.Range("J86").Value = 26101.75
$ws.Range("K86").Value = 1966.6666
$ws.Range("L86").Value = 26101.75
$ws.Range("M86").Value = -843.6666
$ws.Range("N86").Value = -28347.75
$ws.Range("H89").Value = 11620.7
$ws.Range("I89").Value = 1966.6666
$ws.Range("J89").Value = 26101.75
$ws.Range("K89").Value = 9833.333000000001
$ws.Range("L89").Value = 130508.75
$ws.Range("M89").Value = -4217.333000000001
$ws.Range("N89").Value = -141740.75
$ws.Range("H94").Value = 2652.25
$ws.Range("I94").Value = 2203
$ws.Range("K94").Value = 2203
$ws.Range("M94").Value = -1752
$ws.Range("H134").Value = 4051.6365
$ws.Range("I134").Value = 2983.2
$ws.Range("J134").Value = 6341.143
$ws.Range("K134").Value = 8949.599999999999
$ws.Range("L134").Value = 19023.429
$ws.Range("M134").Value = -6414.599999999999
$ws.Range("N134").Value = -24093.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1000000000
$ws.Range("J4").Value = 1000000000
$ws.Range("L4").Value = 1000000000
$ws.Range("N4").Value = -1000000224
$ws.Range("H22").Value = 421.35715
$ws.Range("I22").Value = 178.85715
$ws.Range("J22").Value = 663.8570999999999
$ws.Range("K22").Value = 178.85715
$ws.Range("L22").Value = 663.8570999999999
$ws.Range("M22").Value = 171.14285
$ws.Range("N22").Value = -1363.8571
$ws.Range("H132").Value = 2697.3914
$ws.Range("I132").Value = 2372.8206
$ws.Range("K132").Value = 7118.4618
$ws.Range("M132").Value = -4588.4618

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 414.85715
$ws.Range("J98").Value = 301.33334
$ws.Range("L98").Value = 904.0000200000001
$ws.Range("N98").Value = -3900.00002
$ws.Range("H107").Value = 376.20514
$ws.Range("I107").Value = 393.56
$ws.Range("J107").Value = 345.2143
$ws.Range("K107").Value = 1180.68
$ws.Range("L107").Value = 1035.6429
$ws.Range("M107").Value = 739.3199999999999
$ws.Range("N107").Value = -4875.6429
$ws.Range("H121").Value = 238.5238
$ws.Range("I121").Value = 168.89473
$ws.Range("J121").Value = 900
$ws.Range("K121").Value = 506.6841900000001
$ws.Range("L121").Value = 2700
$ws.Range("M121").Value = 803.3158099999999
$ws.Range("N121").Value = -5320
$ws.Range("H128").Value = 470000
$ws.Range("I128").Value = 470000
$ws.Range("K128").Value = 1410000
$ws.Range("M128").Value = -1405020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1642.9642
$ws.Range("I122").Value = 1270.6818
$ws.Range("J122").Value = 3008
$ws.Range("K122").Value = 3812.0454
$ws.Range("L122").Value = 9024
$ws.Range("M122").Value = -1362.0454
$ws.Range("N122").Value = -13924
$ws.Range("H132").Value = 2687.0386
$ws.Range("I132").Value = 3922.4546
$ws.Range("J132").Value = 1781.0667
$ws.Range("K132").Value = 11767.3638
$ws.Range("L132").Value = 5343.2001
$ws.Range("M132").Value = -9237.363799999999
$ws.Range("N132").Value = -10403.2001
$ws.Range("H136").Value = 11022.75
$ws.Range("J136").Value = 10276.134
$ws.Range("L136").Value = 30828.402
$ws.Range("N136").Value = -35928.402

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 24925.666
$ws.Range("I39").Value = 24925.666
$ws.Range("K39").Value = 24925.666
$ws.Range("M39").Value = -24465.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2288.0312
$ws.Range("I132").Value = 2067.7827
$ws.Range("J132").Value = 2850.889
$ws.Range("K132").Value = 6203.348100000001
$ws.Range("L132").Value = 8552.667000000001
$ws.Range("M132").Value = -3673.348100000001
$ws.Range("N132").Value = -13612.667
